# Update sheet1 with refreshed TPM-derived LR-pair data for Wnt11-Fzd8.
# The Sending cluster "ECs" rows are removed (ligand tool no longer reports
# ECs as a sender for this pair), leaving FAPs/MuSCs senders x ECs/FAPs/MuSCs
# targets, and all numeric TPM-derived metrics are refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 8-10 (previously "ECs" sending cluster, now removed).
$ws.Rows("8:10").Delete()

# Refresh rows 2-7 (Sending cluster / Target cluster labels + all TPM metrics).
# Row 2: FAPs -> ECs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Wnt11"
$ws.Cells.Item(2, 3).Value = "Fzd8"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 20.223983
$ws.Cells.Item(2, 8).Value = 60.671949
$ws.Cells.Item(2, 9).Value = 0.9624502889455165
$ws.Cells.Item(2, 10).Value = 0.9624502889455167
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.235341333333333
$ws.Cells.Item(2, 14).Value = 9.706024
$ws.Cells.Item(2, 15).Value = 0.2153734454473681
$ws.Cells.Item(2, 16).Value = 0.2153734454473681
$ws.Cells.Item(2, 17).Value = 65.43148812453066
$ws.Cells.Item(2, 18).Value = 588.883393120776
$ws.Cells.Item(2, 19).Value = 0.2072862348020109
$ws.Cells.Item(2, 20).Value = 0.2072862348020109

# Row 3: FAPs -> FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Wnt11"
$ws.Cells.Item(3, 3).Value = "Fzd8"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 20.223983
$ws.Cells.Item(3, 8).Value = 60.671949
$ws.Cells.Item(3, 9).Value = 0.9624502889455165
$ws.Cells.Item(3, 10).Value = 0.9624502889455167
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.273511666666667
$ws.Cells.Item(3, 14).Value = 21.820535
$ws.Cells.Item(3, 15).Value = 0.4841904166376352
$ws.Cells.Item(3, 16).Value = 0.4841904166376352
$ws.Cells.Item(3, 17).Value = 147.0993762969683
$ws.Cells.Item(3, 18).Value = 1323.894386672715
$ws.Cells.Item(3, 19).Value = 0.4660092063975421
$ws.Cells.Item(3, 20).Value = 0.4660092063975421

# Row 4: FAPs -> MuSCs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Wnt11"
$ws.Cells.Item(4, 3).Value = "Fzd8"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 20.223983
$ws.Cells.Item(4, 8).Value = 60.671949
$ws.Cells.Item(4, 9).Value = 0.9624502889455165
$ws.Cells.Item(4, 10).Value = 0.9624502889455167
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.513153666666667
$ws.Cells.Item(4, 14).Value = 13.539461
$ws.Cells.Item(4, 15).Value = 0.3004361379149967
$ws.Cells.Item(4, 16).Value = 0.3004361379149967
$ws.Cells.Item(4, 17).Value = 91.27394303105433
$ws.Cells.Item(4, 18).Value = 821.4654872794889
$ws.Cells.Item(4, 19).Value = 0.2891548477459637
$ws.Cells.Item(4, 20).Value = 0.2891548477459637

# Row 5: MuSCs -> ECs
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Wnt11"
$ws.Cells.Item(5, 3).Value = "Fzd8"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.7890326666666668
$ws.Cells.Item(5, 8).Value = 2.367098
$ws.Cells.Item(5, 9).Value = 0.03754971105448342
$ws.Cells.Item(5, 10).Value = 0.03754971105448343
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.235341333333333
$ws.Cells.Item(5, 14).Value = 9.706024
$ws.Cells.Item(5, 15).Value = 0.2153734454473681
$ws.Cells.Item(5, 16).Value = 0.2153734454473681
$ws.Cells.Item(5, 17).Value = 2.552789999816889
$ws.Cells.Item(5, 18).Value = 22.975109998352
$ws.Cells.Item(5, 19).Value = 0.00808721064535722
$ws.Cells.Item(5, 20).Value = 0.008087210645357223

# Row 6: MuSCs -> FAPs
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Wnt11"
$ws.Cells.Item(6, 3).Value = "Fzd8"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.7890326666666668
$ws.Cells.Item(6, 8).Value = 2.367098
$ws.Cells.Item(6, 9).Value = 0.03754971105448342
$ws.Cells.Item(6, 10).Value = 0.03754971105448343
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.273511666666667
$ws.Cells.Item(6, 14).Value = 21.820535
$ws.Cells.Item(6, 15).Value = 0.4841904166376352
$ws.Cells.Item(6, 16).Value = 0.4841904166376352
$ws.Cells.Item(6, 17).Value = 5.739038306381112
$ws.Cells.Item(6, 18).Value = 51.65134475743001
$ws.Cells.Item(6, 19).Value = 0.01818121024009315
$ws.Cells.Item(6, 20).Value = 0.01818121024009315

# Row 7: MuSCs -> MuSCs
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Wnt11"
$ws.Cells.Item(7, 3).Value = "Fzd8"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.7890326666666668
$ws.Cells.Item(7, 8).Value = 2.367098
$ws.Cells.Item(7, 9).Value = 0.03754971105448342
$ws.Cells.Item(7, 10).Value = 0.03754971105448343
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.513153666666667
$ws.Cells.Item(7, 14).Value = 13.539461
$ws.Cells.Item(7, 15).Value = 0.3004361379149967
$ws.Cells.Item(7, 16).Value = 0.3004361379149967
$ws.Cells.Item(7, 17).Value = 3.561025672686445
$ws.Cells.Item(7, 18).Value = 32.049231054178
$ws.Cells.Item(7, 19).Value = 0.01128129016903306
$ws.Cells.Item(7, 20).Value = 0.01128129016903306

$ws.Range("A1").Select()
